$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'CreatedAt: 2025-06-01T18:07:33'
$ws.Range("V4").Value = 71.18000000000001
$ws.Range("W4").Value = 63.31
$ws.Range("X4").Value = 14.81
$ws.Range("Y4").Value = 14.5
$ws.Range("Z4").Value = 6.38
$ws.Range("V6").Value = -0.28
$ws.Range("W6").Value = 0.06
$ws.Range("X6").Value = 0.03
$ws.Range("Z6").Value = 0.1
$ws.Range("V8").Value = 0
$ws.Range("W8").Value = 0
$ws.Range("V9").Value = 71.59999999999999
$ws.Range("W9").Value = 74.34
$ws.Range("X9").Value = 15.36
$ws.Range("Y9").Value = 15.09
$ws.Range("Z9").Value = 6.61
$ws.Range("V11").Value = 0.14
$ws.Range("W11").Value = 2.02
$ws.Range("X11").Value = 0.58
$ws.Range("Y11").Value = 0.78
$ws.Range("Z11").Value = 0.32
$ws.Range("W12").Value = 9.07
$ws.Range("V13").Value = 0
$ws.Range("W13").Value = 0
$ws.Range("V14").Value = 72.67
$ws.Range("W14").Value = 75.27
$ws.Range("X14").Value = 15.49
$ws.Range("Y14").Value = 15.21
$ws.Range("Z14").Value = 6.61
$ws.Range("V15").Value = 1.07
$ws.Range("W15").Value = 0.93
$ws.Range("X15").Value = 0.13
$ws.Range("V16").Value = 0.14
$ws.Range("W16").Value = 2.02
$ws.Range("X16").Value = 0.58
$ws.Range("Y16").Value = 0.78
$ws.Range("Z16").Value = 0.32
$ws.Range("W17").Value = 9.07
$ws.Range("V18").Value = 0
$ws.Range("W18").Value = 0
$ws.Range("V19").Value = 71.68000000000001
$ws.Range("W19").Value = 64.15000000000001
$ws.Range("X19").Value = 14.99
$ws.Range("Y19").Value = 14.72
$ws.Range("Z19").Value = 6.49
$ws.Range("V21").Value = 0.22
$ws.Range("W21").Value = 0.9
$ws.Range("X21").Value = 0.21
$ws.Range("Z21").Value = 0.2
$ws.Range("V23").Value = 0
$ws.Range("W23").Value = 0
$ws.Range("V24").Value = 71.68000000000001
$ws.Range("W24").Value = 64.15000000000001
$ws.Range("X24").Value = 14.99
$ws.Range("Y24").Value = 14.72
$ws.Range("Z24").Value = 6.49
$ws.Range("V26").Value = 0.22
$ws.Range("W26").Value = 0.9
$ws.Range("X26").Value = 0.21
$ws.Range("Z26").Value = 0.2
$ws.Range("V28").Value = 0
$ws.Range("W28").Value = 0
$ws.Range("V29").Value = 72.48
$ws.Range("W29").Value = 65.06999999999999
$ws.Range("X29").Value = 15.2
$ws.Range("Y29").Value = 14.98
$ws.Range("Z29").Value = 6.6
$ws.Range("V31").Value = 1.01
$ws.Range("W31").Value = 1.82
$ws.Range("X31").Value = 0.43
$ws.Range("V33").Value = 0
$ws.Range("W33").Value = 0
$ws.Range("V34").Value = 71.89
$ws.Range("W34").Value = 66
$ws.Range("X34").Value = 15.47
$ws.Range("Y34").Value = 15.18
$ws.Range("Z34").Value = 6.64
$ws.Range("V35").Value = 1.07
$ws.Range("W35").Value = 0.93
$ws.Range("X35").Value = 0.13
$ws.Range("V36").Value = -0.64
$ws.Range("W36").Value = 1.82
$ws.Range("X36").Value = 0.57
$ws.Range("Y36").Value = 0.75
$ws.Range("Z36").Value = 0.35
$ws.Range("V38").Value = 0
$ws.Range("W38").Value = 0
$ws.Range("V39").Value = 71.18000000000001
$ws.Range("W39").Value = 63.31
$ws.Range("X39").Value = 14.81
$ws.Range("Y39").Value = 14.5
$ws.Range("Z39").Value = 6.38
$ws.Range("V41").Value = -0.28
$ws.Range("W41").Value = 0.06
$ws.Range("X41").Value = 0.03
$ws.Range("Z41").Value = 0.1
$ws.Range("V43").Value = 0
$ws.Range("W43").Value = 0
$ws.Range("V44").Value = 72.26000000000001
$ws.Range("W44").Value = 62.68
$ws.Range("X44").Value = 14.69
$ws.Range("Y44").Value = 14.19
$ws.Range("Z44").Value = 6.18
$ws.Range("V46").Value = 0.79
$ws.Range("W46").Value = -0.5600000000000001
$ws.Range("V48").Value = 0
$ws.Range("W48").Value = 0
$ws.Range("V49").Value = 61.6
$ws.Range("W49").Value = 53.97
$ws.Range("X49").Value = 13.23
$ws.Range("Y49").Value = 12.72
$ws.Range("Z49").Value = 5.51
$ws.Range("V51").Value = -9.859999999999999
$ws.Range("W51").Value = -9.279999999999999
$ws.Range("X51").Value = -1.55
$ws.Range("Z51").Value = -0.78
$ws.Range("V53").Value = 0
$ws.Range("W53").Value = 0
$ws.Range("V54").Value = 63.58
$ws.Range("W54").Value = 56.17
$ws.Range("X54").Value = 13.25
$ws.Range("Y54").Value = 12.81
$ws.Range("Z54").Value = 5.69
$ws.Range("V56").Value = -7.88
$ws.Range("W56").Value = -7.08
$ws.Range("X56").Value = -1.52
$ws.Range("Y56").Value = -1.5
$ws.Range("Z56").Value = -0.6
$ws.Range("V58").Value = 0
$ws.Range("W58").Value = 0
$ws.Range("V59").Value = 74.44
$ws.Range("W59").Value = 64.15000000000001
$ws.Range("X59").Value = 15.05
$ws.Range("Y59").Value = 14.51
$ws.Range("Z59").Value = 6.28
$ws.Range("V61").Value = 2.98
$ws.Range("W61").Value = 0.9
$ws.Range("V63").Value = 0
$ws.Range("W63").Value = 0
$ws.Range("V64").Value = 75.54000000000001
$ws.Range("W64").Value = 65
$ws.Range("X64").Value = 15.25
$ws.Range("Y64").Value = 14.67
$ws.Range("Z64").Value = 6.35
$ws.Range("V66").Value = 4.08
$ws.Range("W66").Value = 1.76
$ws.Range("V68").Value = 0
$ws.Range("W68").Value = 0
$ws.Range("V69").Value = 76.43000000000001
$ws.Range("W69").Value = 65.06999999999999
$ws.Range("X69").Value = 15.33
$ws.Range("Y69").Value = 14.75
$ws.Range("Z69").Value = 6.33
$ws.Range("V71").Value = 4.97
$ws.Range("W71").Value = 1.82
$ws.Range("V73").Value = 0
$ws.Range("W73").Value = 0
$ws.Range("V74").Value = 73.06999999999999
$ws.Range("W74").Value = 63.06
$ws.Range("X74").Value = 14.79
$ws.Range("Y74").Value = 14.25
$ws.Range("Z74").Value = 6.17
$ws.Range("V76").Value = 1.61
$ws.Range("W76").Value = -0.19
$ws.Range("V78").Value = 0
$ws.Range("W78").Value = 0
$ws.Range("V79").Value = 71.45999999999999
$ws.Range("W79").Value = 63.25
$ws.Range("X79").Value = 14.78
$ws.Range("Y79").Value = 14.31
$ws.Range("Z79").Value = 6.29
$ws.Range("V83").Value = 0
$ws.Range("W83").Value = 0
$ws.Range("V84").Value = 62.85
$ws.Range("W84").Value = 55.63
$ws.Range("X84").Value = 13.18
$ws.Range("Y84").Value = 12.73
$ws.Range("Z84").Value = 5.67
$ws.Range("V86").Value = -8.609999999999999
$ws.Range("W86").Value = -7.62
$ws.Range("X86").Value = -1.6
$ws.Range("Y86").Value = -1.58
$ws.Range("Z86").Value = -0.62
$ws.Range("V88").Value = 0
$ws.Range("W88").Value = 0
$ws.Range("V89").Value = 72.48
$ws.Range("W89").Value = 65.06999999999999
$ws.Range("X89").Value = 15.2
$ws.Range("Y89").Value = 14.98
$ws.Range("Z89").Value = 6.6
$ws.Range("V91").Value = 1.01
$ws.Range("W91").Value = 1.82
$ws.Range("X91").Value = 0.43
$ws.Range("V93").Value = 0
$ws.Range("W93").Value = 0
